$d = $word.ActiveDocument

# Paragraph 79 (paraId 79918D13) currently sits right before the final
# bookmark paragraph (paraId 63474020, index 80). It is empty with style
# "Note Level 1" and an explicit (no-op) numPr (ilvl=0, numId=0). We add
# the "Eager Learner" text to it, keeping its existing paragraph
# formatting intact.
$p = $d.Paragraphs.Item(79)
$p.Range.InsertBefore("Eager Learner – once model is created you can discard the data.")

# Insert "Lazy Learner" as a new paragraph after it, same style/numPr.
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item(80)
$p.Range.InsertBefore("Lazy Learner – can’t discard the data")

# Empty paragraph, same style/numPr.
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item(81)

# "Batch Learner / Inline or Incremental Learner", same style/numPr.
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item(82)
$p.Range.InsertBefore("Batch Learner / Inline or Incremental Learner")

# Empty paragraph, same style/numPr.
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item(83)

# "* know the biases of each algorithm" -- plain "Note Level 1" (no numPr).
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item(84)
$p.Style = "Note Level 1"
$p.Range.InsertBefore("* know the biases of each algorithm")

# Empty paragraph, plain "Note Level 1".
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item(85)

# "K-Nearest Neighbor:", plain "Note Level 1".
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item(86)
$p.Range.InsertBefore("K-Nearest Neighbor:")

# "Bias: things that look like me should be classified like me." -- "Note Level 2".
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item(87)
$p.Style = "Note Level 2"
$p.Range.InsertBefore("Bias: things that look like me should be classified like me.")

# "Assumption : we can put things into some meaningful space and the
# distance around them is meaningful.  " -- "Note Level 2".
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item(88)
$p.Range.InsertBefore("Assumption : we can put things into some meaningful space and the distance around them is meaningful.  ")

# "Assumptoin : we have continuious space" -- "Note Level 2".
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item(89)
$p.Range.InsertBefore("Assumptoin : we have continuious space")

# "….What is the distance between red and green?" -- "Note Level 2".
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item(90)
$p.Range.InsertBefore("….What is the distance between red and green?")

# "How do you measure these things." -- "Note Level 3".
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item(91)
$p.Style = "Note Level 3"
$p.Range.InsertBefore("How do you measure these things.")

# "Could do ordinal – discrete but the difference makes sense." -- "Note Level 3".
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item(92)
$p.Range.InsertBefore("Could do ordinal – discrete but the difference makes sense.")

# Finally, the paragraph that holds the _GoBack bookmark (was index 80,
# now pushed down to 93) becomes "Note Level 2" and loses its (no-op)
# numPr list formatting.
$p = $d.Paragraphs.Item(93)
$p.Style = "Note Level 2"
